$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Add the new row of data for 30 April 2020 (row 50)
$ws.Range("A50").Value = 43951
$ws.Range("B50").Value = 42004
$ws.Range("C50").Value = 2615
$ws.Range("D50").Value = 93
$ws.Range("E50").Value = 4846

# Grow the worksheet Table ("Table3") to include the new row
$lo = $ws.ListObjects("Table3")
$lo.Resize($ws.Range("A1:E50"))

# Move the active selection to match post-edit state
$ws.Range("B51").Select()
